# Update Sheets via scheduled runner
# Applies updated profit/price figures to several worksheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 80462.57000000001
$ws.Range("I98").Value = 80462.57000000001
$ws.Range("K98").Value = 80462.57000000001
$ws.Range("M98").Value = -78964.57000000001

$ws.Range("H116").Value = 3825
$ws.Range("I116").Value = 1666.6666
$ws.Range("J116").Value = 5120
$ws.Range("K116").Value = 1666.6666
$ws.Range("L116").Value = 5120
$ws.Range("M116").Value = 1775.3334
$ws.Range("N116").Value = -12004

$ws.Range("H122").Value = 80462.57000000001
$ws.Range("I122").Value = 80462.57000000001
$ws.Range("K122").Value = 241387.71
$ws.Range("M122").Value = -238937.71

$ws.Range("H137").Value = 1668.2084
$ws.Range("I137").Value = 1310.3334
$ws.Range("J137").Value = 2026.0834
$ws.Range("K137").Value = 3931.0002
$ws.Range("L137").Value = 6078.2502
$ws.Range("M137").Value = -1381.0002
$ws.Range("N137").Value = -11178.2502

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 71429784
$ws.Range("J2").Value = 1717.5
$ws.Range("L2").Value = 1717.5
$ws.Range("N2").Value = -1943.5

$ws.Range("H32").Value = 15712.862
$ws.Range("I32").Value = 15559.75
$ws.Range("K32").Value = 15559.75
$ws.Range("M32").Value = -15272.75

$ws.Range("H45").Value = 41667516
$ws.Range("I45").Value = 55556240
$ws.Range("J45").Value = 1350
$ws.Range("K45").Value = 55556240
$ws.Range("L45").Value = 1350
$ws.Range("M45").Value = -55555863
$ws.Range("N45").Value = -2104

$ws.Range("H61").Value = 2011.24
$ws.Range("I61").Value = 1666.5
$ws.Range("K61").Value = 1666.5
$ws.Range("M61").Value = -1454.5

$ws.Range("H116").Value = 71429784
$ws.Range("J116").Value = 1717.5
$ws.Range("L116").Value = 1717.5
$ws.Range("N116").Value = -6305.5

$ws.Range("H122").Value = 1986.8182
$ws.Range("I122").Value = 1886.375
$ws.Range("J122").Value = 2254.6667
$ws.Range("K122").Value = 5659.125
$ws.Range("L122").Value = 6764.000100000001
$ws.Range("M122").Value = -3209.125
$ws.Range("N122").Value = -11664.0001

$ws.Range("H132").Value = 4055.3584
$ws.Range("I132").Value = 4266.8716
$ws.Range("J132").Value = 3466.1428
$ws.Range("K132").Value = 12800.6148
$ws.Range("L132").Value = 10398.4284
$ws.Range("M132").Value = -10270.6148
$ws.Range("N132").Value = -15458.4284

$ws.Range("H136").Value = 2011.24
$ws.Range("I136").Value = 1666.5
$ws.Range("K136").Value = 4999.5
$ws.Range("M136").Value = -2449.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 71429784
$ws.Range("J3").Value = 1717.5
$ws.Range("L3").Value = 1717.5
$ws.Range("N3").Value = -1945.5

$ws.Range("H20").Value = 4389.2646
$ws.Range("I20").Value = 4884.5
$ws.Range("J20").Value = 3200.7
$ws.Range("K20").Value = 4884.5
$ws.Range("L20").Value = 3200.7
$ws.Range("M20").Value = -4637.5
$ws.Range("N20").Value = -3694.7

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4390085.5
$ws.Range("I31").Value = 3304.75
$ws.Range("J31").Value = 11910282
$ws.Range("K31").Value = 3304.75
$ws.Range("L31").Value = 11910282
$ws.Range("M31").Value = -3009.75
$ws.Range("N31").Value = -11910872

$ws.Range("H34").Value = 4390085.5
$ws.Range("I34").Value = 3304.75
$ws.Range("J34").Value = 11910282
$ws.Range("K34").Value = 3304.75
$ws.Range("L34").Value = 11910282
$ws.Range("M34").Value = -3102.75
$ws.Range("N34").Value = -11910686

$ws.Range("H109").Value = 10000
$ws.Range("J109").Value = 10000
$ws.Range("L109").Value = 10000
$ws.Range("N109").Value = -12080

$ws.Range("H132").Value = 5003348
$ws.Range("I132").Value = 3077
$ws.Range("J132").Value = 8336861.5
$ws.Range("K132").Value = 9231
$ws.Range("L132").Value = 25010584.5
$ws.Range("M132").Value = -6701
$ws.Range("N132").Value = -25015644.5

$ws.Range("H134").Value = 1266
$ws.Range("I134").Value = 1184
$ws.Range("J134").Value = 1375.3334
$ws.Range("K134").Value = 3552
$ws.Range("L134").Value = 4126.0002
$ws.Range("M134").Value = -1017
$ws.Range("N134").Value = -9196.0002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 741.26
$ws.Range("I131").Value = 315.45
$ws.Range("J131").Value = 847.7125
$ws.Range("K131").Value = 946.3499999999999
$ws.Range("L131").Value = 2543.1375
$ws.Range("M131").Value = 4093.65
$ws.Range("N131").Value = -12623.1375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2065.4443
$ws.Range("I97").Value = 1676.6842
$ws.Range("J97").Value = 2988.75
$ws.Range("K97").Value = 1676.6842
$ws.Range("L97").Value = 2988.75
$ws.Range("M97").Value = -1180.6842
$ws.Range("N97").Value = -3980.75

$ws.Range("H122").Value = 76926290
$ws.Range("I122").Value = 166670460
$ws.Range("J122").Value = 2714
$ws.Range("K122").Value = 500011380
$ws.Range("L122").Value = 8142
$ws.Range("M122").Value = -500008930
$ws.Range("N122").Value = -13042

$ws.Range("H132").Value = 2692.0557
$ws.Range("I132").Value = 2230.6667
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 6692.000100000001
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -4162.000100000001
$ws.Range("N132").Value = -20057

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 9808.115
$ws.Range("I132").Value = 18039.834
$ws.Range("J132").Value = 2752.3572
$ws.Range("K132").Value = 54119.50199999999
$ws.Range("L132").Value = 8257.071599999999
$ws.Range("M132").Value = -51589.50199999999
$ws.Range("N132").Value = -13317.0716

$ws.Range("H136").Value = 4724.2705
$ws.Range("I136").Value = 5596.423
$ws.Range("J136").Value = 2662.818
$ws.Range("K136").Value = 16789.269
$ws.Range("L136").Value = 7988.454000000001
$ws.Range("M136").Value = -14239.269
$ws.Range("N136").Value = -13088.454

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 456.25
$ws.Range("I113").Value = 455.57144
$ws.Range("J113").Value = 457.2
$ws.Range("K113").Value = 1366.71432
$ws.Range("L113").Value = 1371.6
$ws.Range("M113").Value = 803.28568
$ws.Range("N113").Value = -5711.6

$ws.Range("H138").Value = 64607.25
$ws.Range("J138").Value = 64607.25
$ws.Range("L138").Value = 64607.25
$ws.Range("N138").Value = -74887.25

